$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Swap columns B and C (both their contents and their widths).
#    Using Cut + Insert moves the whole column (values, shared-string
#    refs and width) as a unit, which keeps the original column width
#    values byte-exact instead of going through a lossy re-computation.
# ------------------------------------------------------------------
$ws.Columns.Item(3).Cut()
$ws.Columns.Item(2).Insert()
$excel.CutCopyMode = $false

# ------------------------------------------------------------------
# 2. New header row (D1:J1) - added in this order so new shared
#    strings line up with how the workbook was actually authored.
# ------------------------------------------------------------------
$ws.Range("D1").Value = "UPC"
$ws.Range("E1").Value = "SECONDARY_UPC"
$ws.Range("F1").Value = "BRAND_NAME"
$ws.Range("G1").Value = "MODEL_NO"
$ws.Range("H1").Value = "MANUFACTURER"
$ws.Range("I1").Value = "UNIT_COST"
$ws.Range("J1").Value = "ORIGINAL_RETAIL"

# New product info rows for the two SKUs already in the sheet.
$ws.Range("F2").Value = "Alfalfa"
$ws.Range("H2").Value = "Asus"

# Last header cell.
$ws.Range("K1").Value = "ENT_STREET_DATE"

# ------------------------------------------------------------------
# 3. Remaining data cells.
# ------------------------------------------------------------------
$ws.Range("D2").Value = 12345678910
$ws.Range("D3").Value = 12345678912

$ws.Range("F3").Value = "Alfalfa"

$ws.Range("G2").Value = 1032423
$ws.Range("G3").Value = 1032423

$ws.Range("H3").Value = "Asus"

$ws.Range("I2").Value = 1.26
$ws.Range("I3").Value = 1.26

$ws.Range("J2").Value = 2.89
$ws.Range("J3").Value = 2.89

# Street-date column, formatted as a date (numFmtId 14 / "mm-dd-yy").
# K3 picks up its format from K2 via copy/paste-special so both cells
# share a single cell style entry instead of allocating two.
$ws.Range("K2").Value = 44677
$ws.Range("K2").NumberFormat = "mm-dd-yy"
$ws.Range("K2").Copy()
$ws.Range("K3").PasteSpecial(-4122)
$ws.Range("K3").Value = 45443
$excel.CutCopyMode = $false

# ------------------------------------------------------------------
# 4. Column widths for the newly added columns (best-fit-style sizing).
# ------------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 11.33
$ws.Columns.Item(5).ColumnWidth = 15.17
$ws.Columns.Item(6).ColumnWidth = 12
$ws.Columns.Item(7).ColumnWidth = 9.83
$ws.Columns.Item(8).ColumnWidth = 14
$ws.Columns.Item(9).ColumnWidth = 9.5

# ------------------------------------------------------------------
# 5. Selection shown in the diff.
# ------------------------------------------------------------------
$ws.Range("G23").Select()
